$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 553 entirely (the "あなたの声を聞かせてください…" post entry),
# which shifts all subsequent rows up by one.
$ws.Rows.Item(553).Delete()
